# ---------------------------------------------------------------------------
# Adds a new "2022-Q4" sheet (holding fund data) right after "总计", shifts
# all the quarterly sheets that follow it back by one tab, and updates the
# summary ("总计") sheet with a new row of totals for 2022-Q4.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force the cell to keep a literal/text value (e.g. fund codes with
    # leading zeros, or numeric-looking names) instead of Excel's default
    # auto-coercion to a Number when a numeric-looking string is assigned.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right before "2022-Q3".
# ---------------------------------------------------------------------------

$beforeSheet = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Add($beforeSheet)
$q4Sheet.Name = "2022-Q4"

# Match page margins used by the other quarterly sheets (values are in
# points: 0.75in=54, 1in=72, 0.5in=36).
$q4Sheet.PageSetup.LeftMargin = 54
$q4Sheet.PageSetup.RightMargin = 54
$q4Sheet.PageSetup.TopMargin = 72
$q4Sheet.PageSetup.BottomMargin = 72
$q4Sheet.PageSetup.HeaderMargin = 36
$q4Sheet.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 2. Populate the new "2022-Q4" sheet with its header + fund rows, copying
#    the header style ("s=2") from an already-existing quarterly sheet.
# ---------------------------------------------------------------------------

$q4Header = @(
    "基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名"
)

$q4Rows = @(
    @("010947", "中欧嘉选混合A", "12.93", "83.49", "2.95", "0.3814", "8"),
    @("159667", "国泰中证机床ETF", "3.49", "99.21", "6.20", "0.2164", "2"),
    @("159663", "华夏中证机床ETF", "1.01", "97.42", "6.12", "0.0618", "2"),
    @("159743", "博时中证湖北新旧动能转换ETF", "3.43", "99.18", "1.63", "0.0559", "10"),
    @("010948", "中欧嘉选混合C", "0.74", "83.49", "2.95", "0.0218", "8")
)

# Copy the bold/bordered header format ("s=2") onto row 1, columns B:H.
$wb.Worksheets.Item("2022-Q3").Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)

for ($col = 2; $col -le 8; $col++) {
    $cell = $q4Sheet.Cells.Item(1, $col)
    Set-TextValue $cell $q4Header[$col - 2]
}

$r = 2
foreach ($row in $q4Rows) {
    # Column A: plain 0-based numeric index.
    $q4Sheet.Cells.Item($r, 1).Value = $r - 2

    # Columns B:G: text values (fund code / name / scale / position / etc).
    for ($col = 2; $col -le 7; $col++) {
        $cell = $q4Sheet.Cells.Item($r, $col)
        Set-TextValue $cell $row[$col - 2]
    }

    # Column H: numeric rank.
    $q4Sheet.Cells.Item($r, 8).Value = [double]$row[6]

    $r++
}

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: insert a new row for 2022-Q4 right
#    after the header, and renumber the leading index column (0-based).
# ---------------------------------------------------------------------------

$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# Re-apply the style ("s=2") that belongs on column A of every data row.
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 5
$summary.Cells.Item(2, 4).Value = 0.74

# Renumber column A (0-based row index) for every row that follows, since
# the "2020-Q4" row is now one row further down (row 10 instead of row 9).
for ($row = 3; $row -le 10; $row++) {
    $summary.Cells.Item($row, 1).Value = $row - 2
}

# ---------------------------------------------------------------------------
# 4. Keep "2020-Q4" as the selected/active tab, like in the source file.
# ---------------------------------------------------------------------------

$wb.Worksheets.Item("2020-Q4").Activate()

Write-Output "2022-Q4 sheet inserted; summary sheet updated."
